$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 992
$ws.Range("F8").Value = 563
$ws.Range("F9").Value = 1467
$ws.Range("F11").Value = 1362
$ws.Range("F12").Value = 3019
$ws.Range("F13").Value = 463
$ws.Range("F14").Value = 1650
$ws.Range("F15").Value = 1366
$ws.Range("F16").Value = 806
$ws.Range("F17").Value = 243
$ws.Range("F18").Value = 1401
$ws.Range("F21").Value = 1132
$ws.Range("F23").Value = 405
$ws.Range("F24").Value = 11
$ws.Range("F25").Value = 3521
$ws.Range("F28").Value = 1561

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 171
$ws.Range("F6").Value = 13
$ws.Range("F7").Value = 55
$ws.Range("F8").Value = 28
$ws.Range("F9").Value = 22
$ws.Range("F12").Value = 83
$ws.Range("F13").Value = 18

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 13

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 13
$ws.Range("F8").Value = 171
$ws.Range("F10").Value = 13
$ws.Range("F11").Value = 55
$ws.Range("F13").Value = 28
$ws.Range("F14").Value = 22
$ws.Range("F15").Value = 992
$ws.Range("F18").Value = 563
$ws.Range("F19").Value = 1467
$ws.Range("F21").Value = 1362
$ws.Range("F22").Value = 3019
$ws.Range("F23").Value = 463
$ws.Range("F24").Value = 1650
$ws.Range("F25").Value = 1366
$ws.Range("F26").Value = 806
$ws.Range("F27").Value = 243
$ws.Range("F28").Value = 1401
$ws.Range("F33").Value = 1132
$ws.Range("F35").Value = 405
$ws.Range("F36").Value = 11
$ws.Range("F37").Value = 3522
$ws.Range("F40").Value = 1561
$ws.Range("F41").Value = 83
$ws.Range("F42").Value = 18
